$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 84, shifting existing rows 84..211 down to 85..212
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new data point
$ws.Cells.Item(84, 1).Value = 10
$ws.Cells.Item(84, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(84, 3).Value = "La Araucanía"
$ws.Cells.Item(84, 4).Value = 44721
$ws.Cells.Item(84, 5).Value = 9
$ws.Cells.Item(84, 6).Value = 100112005
$ws.Cells.Item(84, 7).Value = "Puerro"
$ws.Cells.Item(84, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 30
$ws.Cells.Item(84, 11).Value = 16000
$ws.Cells.Item(84, 12).Value = 16000
$ws.Cells.Item(84, 13).Value = 16000
$ws.Cells.Item(84, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(84, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(84, 16).Value = 1333
$ws.Cells.Item(84, 17).Value = 12
$ws.Cells.Item(84, 18).Value = "Hortaliza"
